# Apply the daily cryptos-list refresh: updates Price (D) and Volume(1h) (E)
# columns for most rows, and for rows 39/40 the two coins (FraxShare and
# InternetComputer(DFINITY)) swapped list positions, so B/C/D/E all change there.
#
# Many "Price" values look numeric (e.g. "1.001") but must stay plain text,
# matching how the source data already stores them (t="inlineStr" in the xlsx).
# Assigning such a string straight to .Value lets Excel auto-convert it to a
# Number, so each text-look-alike is written via a brief NumberFormat="@"
# (Text) override, then ClearFormats() drops that temporary formatting again
# (cells here carry no real formatting of their own) while the stored cell
# stays Text-typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "28.183.44"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.26%  "

Set-TextValue $ws.Cells.Item(3, 4) "1.797.46"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.52%  "

Set-TextValue $ws.Cells.Item(4, 4) "1.001"
Set-TextValue $ws.Cells.Item(4, 5) "  -0.07%  "

Set-TextValue $ws.Cells.Item(5, 4) "314.43"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.41%  "

Set-TextValue $ws.Cells.Item(6, 4) "1.001"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.04%  "

Set-TextValue $ws.Cells.Item(7, 4) "0.5211"
Set-TextValue $ws.Cells.Item(7, 5) "  +2.01%  "

Set-TextValue $ws.Cells.Item(8, 4) "0.3822"
Set-TextValue $ws.Cells.Item(8, 5) "  -3.26%  "

Set-TextValue $ws.Cells.Item(9, 4) "0.07937"
Set-TextValue $ws.Cells.Item(9, 5) "  -3.62%  "

Set-TextValue $ws.Cells.Item(10, 4) "41.41"
Set-TextValue $ws.Cells.Item(10, 5) "  -0.91%  "

Set-TextValue $ws.Cells.Item(11, 5) "  -1.29%  "

Set-TextValue $ws.Cells.Item(12, 4) "6.288"
Set-TextValue $ws.Cells.Item(12, 5) "  -1.16%  "

Set-TextValue $ws.Cells.Item(13, 5) "  -0.02%  "

Set-TextValue $ws.Cells.Item(14, 4) "20.57"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.79%  "

Set-TextValue $ws.Cells.Item(15, 4) "1.792.05"
Set-TextValue $ws.Cells.Item(15, 5) "  -1.46%  "

Set-TextValue $ws.Cells.Item(16, 4) "7.257"
Set-TextValue $ws.Cells.Item(16, 5) "  -4.08%  "

Set-TextValue $ws.Cells.Item(17, 4) "93.29"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.31%  "

Set-TextValue $ws.Cells.Item(18, 4) "0.00001085"
Set-TextValue $ws.Cells.Item(18, 5) "  -3.44%  "

Set-TextValue $ws.Cells.Item(19, 4) "0.06558"
Set-TextValue $ws.Cells.Item(19, 5) "  -1.57%  "

Set-TextValue $ws.Cells.Item(20, 4) "1.001"
Set-TextValue $ws.Cells.Item(20, 5) "  +0.01%  "

Set-TextValue $ws.Cells.Item(21, 4) "17.30"
Set-TextValue $ws.Cells.Item(21, 5) "  -3.03%  "

Set-TextValue $ws.Cells.Item(22, 4) "5.953"
Set-TextValue $ws.Cells.Item(22, 5) "  -2.52%  "

Set-TextValue $ws.Cells.Item(23, 4) "28.228.38"
Set-TextValue $ws.Cells.Item(23, 5) "  -1.25%  "

Set-TextValue $ws.Cells.Item(24, 4) "11.15"
Set-TextValue $ws.Cells.Item(24, 5) "  -2.46%  "

Set-TextValue $ws.Cells.Item(25, 4) "2.268"
Set-TextValue $ws.Cells.Item(25, 5) "  -0.18%  "

Set-TextValue $ws.Cells.Item(26, 4) "160.95"
Set-TextValue $ws.Cells.Item(26, 5) "  +2.45%  "

Set-TextValue $ws.Cells.Item(27, 4) "20.46"
Set-TextValue $ws.Cells.Item(27, 5) "  -4.23%  "

Set-TextValue $ws.Cells.Item(28, 4) "2.000.52"
Set-TextValue $ws.Cells.Item(28, 5) "  -1.44%  "

Set-TextValue $ws.Cells.Item(29, 4) "2.336"
Set-TextValue $ws.Cells.Item(29, 5) "  -3.19%  "

Set-TextValue $ws.Cells.Item(30, 4) "123.25"
Set-TextValue $ws.Cells.Item(30, 5) "  -2.75%  "

Set-TextValue $ws.Cells.Item(31, 4) "0.1068"
Set-TextValue $ws.Cells.Item(31, 5) "  -1.80%  "

Set-TextValue $ws.Cells.Item(32, 4) "1.054"
Set-TextValue $ws.Cells.Item(32, 5) "  -5.68%  "

Set-TextValue $ws.Cells.Item(33, 4) "3.672"
Set-TextValue $ws.Cells.Item(33, 5) "  +0.28%  "

Set-TextValue $ws.Cells.Item(34, 4) "5.567"
Set-TextValue $ws.Cells.Item(34, 5) "  -3.79%  "

Set-TextValue $ws.Cells.Item(35, 4) "0.07311"
Set-TextValue $ws.Cells.Item(35, 5) "  +3.33%  "

Set-TextValue $ws.Cells.Item(36, 4) "12.30"
Set-TextValue $ws.Cells.Item(36, 5) "  +8.99%  "

Set-TextValue $ws.Cells.Item(37, 4) "0.02324"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.41%  "

Set-TextValue $ws.Cells.Item(41, 4) "0.6157"
Set-TextValue $ws.Cells.Item(41, 5) "  -3.16%  "

Set-TextValue $ws.Cells.Item(42, 4) "1.164"
Set-TextValue $ws.Cells.Item(42, 5) "  -1.62%  "

Set-TextValue $ws.Cells.Item(43, 4) "1.370"
Set-TextValue $ws.Cells.Item(43, 5) "  -2.28%  "

Set-TextValue $ws.Cells.Item(44, 4) "13.20"
Set-TextValue $ws.Cells.Item(44, 5) "  -3.31%  "

Set-TextValue $ws.Cells.Item(45, 4) "3.782"
Set-TextValue $ws.Cells.Item(45, 5) "  +1.15%  "

Set-TextValue $ws.Cells.Item(46, 4) "0.5999"
Set-TextValue $ws.Cells.Item(46, 5) "  +0.37%  "

Set-TextValue $ws.Cells.Item(47, 4) "127.63"
Set-TextValue $ws.Cells.Item(47, 5) "  +1.86%  "

Set-TextValue $ws.Cells.Item(48, 4) "1.231"
Set-TextValue $ws.Cells.Item(48, 5) "  +2.90%  "

Set-TextValue $ws.Cells.Item(49, 4) "1.920"
Set-TextValue $ws.Cells.Item(49, 5) "  -3.73%  "

Set-TextValue $ws.Cells.Item(50, 4) "0.06772"
Set-TextValue $ws.Cells.Item(50, 5) "  -2.50%  "

Set-TextValue $ws.Cells.Item(51, 4) "73.19"
Set-TextValue $ws.Cells.Item(51, 5) "  -1.57%  "

# Rows 39/40: FraxShare and InternetComputer(DFINITY) swapped list positions
# (row 39 was FraxShare, now InternetComputer(DFINITY); row 40 vice versa),
# with refreshed Price/Volume(1h) values for each.
$ws.Cells.Item(39, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Cells.Item(39, 4) "5.071"
Set-TextValue $ws.Cells.Item(39, 5) "  -3.48%  "

$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Cells.Item(40, 4) "8.622"
Set-TextValue $ws.Cells.Item(40, 5) "  -1.78%  "
